# "working in risk calculator"
# Add two new configuration parameters (ris_years / ris_types_analysis) to the
# "conf" sheet, right after the existing rows, and make the "conf" sheet the
# active/selected one (it was "buffer" before).

$wb     = $excel.ActiveWorkbook
$wsConf = $wb.Worksheets.Item("conf")

# --- new parameter rows on "conf" --------------------------------------
# Row 31: ris_years | 2019 | 2010, 2012, 2019
$wsConf.Range("A31").Value2 = "ris_years"

# Row 32: ris_types_analysis | detail,summary
$wsConf.Range("A32").Value2 = "ris_types_analysis"
$wsConf.Range("B32").Value2 = "detail,summary"

# C31 carries the list of available years, formatted with a number format
# (leftover from pasting into a numeric-formatted cell).
$wsConf.Range("C31").Value2 = "2010, 2012, 2019"
$wsConf.Range("C31").NumberFormat = "#,##0"

# B31 carries the selected year, stored as a number but displayed with a
# text number format (also leftover paste formatting).
$wsConf.Range("B31").Value2 = 2019
$wsConf.Range("B31").NumberFormat = "@"

# --- view / selection state ---------------------------------------------
# Scroll roughly to where the new rows are and leave the selection on B31.
[void]$wsConf.Range("A25").Select()
[void]$wsConf.Range("B31").Select()

# "conf" becomes the active sheet/tab (previously "buffer" was active).
[void]$wsConf.Activate()
